$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry describes the new Price (column D) and/or Volume(1h) (column E)
# value for a given row, as captured from the source diff.
$updates = @(
    @{Row=2; D='26.688.91'; E='  +0.22%  '},
    @{Row=3; D='1.596.66'; E='  +0.70%  '},
    @{Row=4; E='  +0.13%  '},
    @{Row=5; D='211.35'; E='  +0.20%  '},
    @{Row=6; E='  +1.33%  '},
    @{Row=7; E='  +0.17%  '},
    @{Row=8; E='  -0.18%  '},
    @{Row=9; D='0.245'; E='  -1.39%  '},
    @{Row=10; D='19.53'; E='  -0.47%  '},
    @{Row=11; D='0.0842'; E='  +1.06%  '},
    @{Row=12; D='1.821.00'; E='  +0.72%  '},
    @{Row=13; D='1.600.70'; E='  +0.94%  '},
    @{Row=15; E='  -1.31%  '},
    @{Row=16; D='64.77'; E='  +0.26%  '},
    @{Row=17; D='26.683.29'; E='  +0.21%  '},
    @{Row=18; D='0.0₃0728'; E='  -0.22%  '},
    @{Row=19; D='208.16'; E='  +0.04%  '},
    @{Row=20; E='  +0.12%  '},
    @{Row=21; E='  +0.53%  '},
    @{Row=22; E='  -0.12%  '},
    @{Row=23; D='2.35'; E='  -0.83%  '},
    @{Row=24; D='8.86'; E='  -0.44%  '},
    @{Row=25; D='145.44'; E='  -0.94%  '},
    @{Row=26; E='  +0.19%  '},
    @{Row=27; D='7.21'; E='  -2.62%  '},
    @{Row=28; E='  +1.18%  '},
    @{Row=29; D='15.25'; E='  -0.41%  '},
    @{Row=30; E='  +0.46%  '},
    @{Row=31; E='  +0.21%  '},
    @{Row=32; D='3.22'; E='  -0.82%  '},
    @{Row=33; D='0.660'; E='  -3.35%  '},
    @{Row=34; D='2.93'; E='  +0.17%  '},
    @{Row=35; D='1.282.88'; E='  -3.73%  '},
    @{Row=36; D='2.46'; E='  +0.77%  '},
    @{Row=37; E='  -1.11%  '},
    @{Row=38; E='  -0.66%  '},
    @{Row=39; D='0.842'; E='  +1.99%  '},
    @{Row=40; E='  +0.11%  '},
    @{Row=41; D='5.42'; E='  +1.15%  '},
    @{Row=42; E='  +1.33%  '},
    @{Row=43; D='0.787'; E='  +0.39%  '},
    @{Row=44; D='63.52'},
    @{Row=45; D='1.733.37'; E='  +0.71%  '},
    @{Row=46; D='0.905'; E='  +9.29%  '},
    @{Row=47; D='90.05'; E='  +0.32%  '},
    @{Row=48; D='1.60'; E='  -0.74%  '},
    @{Row=49; E='  +2.15%  '},
    @{Row=50; E='  -0.24%  '},
    @{Row=51; D='7.45'; E='  -0.67%  '}
)

foreach ($u in $updates) {
    if ($u.ContainsKey('D')) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($u.ContainsKey('E')) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
